$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 73
$srcRow = 72

# Copy date-format/style from A72 to A73, then set value to next day (45629)
$ws.Cells.Item($srcRow, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item($newRow, 1).Value = 45629

# Fill remaining columns (B..J) with the same values as row 72
$ws.Cells.Item($newRow, 2).Value  = 116.4121952
$ws.Cells.Item($newRow, 3).Value  = 0.00170247
$ws.Cells.Item($newRow, 4).Value  = 0.008850780000000001
$ws.Cells.Item($newRow, 5).Value  = 0.06933635
$ws.Cells.Item($newRow, 6).Value  = 12792.90181321
$ws.Cells.Item($newRow, 7).Value  = 465.80531254
$ws.Cells.Item($newRow, 8).Value  = 0.24
$ws.Cells.Item($newRow, 9).Value  = 1.7904431
$ws.Cells.Item($newRow, 10).Value = 485.38834923
